# Auto-generated edit script applying the Leviathan_Profits market-data refresh diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 816
$ws.Range("I15").Value = 816
$ws.Range("K15").Value = 2448
$ws.Range("M15").Value = -2279
$ws.Range("H17").Value = 910847.4399999999
$ws.Range("J17").Value = 910847.4399999999
$ws.Range("L17").Value = 2732542.32
$ws.Range("N17").Value = -2732878.32
$ws.Range("H55").Value = 246.95238
$ws.Range("I55").Value = 248.33333
$ws.Range("K55").Value = 248.33333
$ws.Range("M55").Value = -34.33332999999999
$ws.Range("H70").Value = 1885.1428
$ws.Range("I70").Value = 1449
$ws.Range("K70").Value = 4347
$ws.Range("M70").Value = -4077
$ws.Range("H73").Value = 1885.1428
$ws.Range("I73").Value = 1449
$ws.Range("K73").Value = 4347
$ws.Range("M73").Value = -3411
$ws.Range("H135").Value = 77639.84
$ws.Range("I135").Value = 608.3
$ws.Range("J135").Value = 334411.66
$ws.Range("K135").Value = 5474.7
$ws.Range("L135").Value = 3009704.94
$ws.Range("M135").Value = -2939.7
$ws.Range("N135").Value = -3014774.94
$ws.Range("H137").Value = 1819.7646
$ws.Range("I137").Value = 1495.2307
$ws.Range("K137").Value = 4485.6921
$ws.Range("M137").Value = -1935.6921
$ws.Range("H138").Value = 2708.1
$ws.Range("I138").Value = 5059
$ws.Range("J138").Value = 2120.375
$ws.Range("K138").Value = 15177
$ws.Range("L138").Value = 6361.125
$ws.Range("M138").Value = -10037
$ws.Range("N138").Value = -16641.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 93308
$ws.Range("I32").Value = 19172.334
$ws.Range("K32").Value = 19172.334
$ws.Range("M32").Value = -18885.334
$ws.Range("H45").Value = 6924.875
$ws.Range("I45").Value = 10386.083
$ws.Range("K45").Value = 10386.083
$ws.Range("M45").Value = -10009.083
$ws.Range("H59").Value = 24995
$ws.Range("J59").Value = 24995
$ws.Range("L59").Value = 24995
$ws.Range("N59").Value = -26603
$ws.Range("H74").Value = 1537.6444
$ws.Range("I74").Value = 1447.8948
$ws.Range("J74").Value = 2024.8572
$ws.Range("K74").Value = 1447.8948
$ws.Range("L74").Value = 2024.8572
$ws.Range("M74").Value = -573.8948
$ws.Range("N74").Value = -3772.8572
$ws.Range("H77").Value = 1537.6444
$ws.Range("I77").Value = 1447.8948
$ws.Range("J77").Value = 2024.8572
$ws.Range("K77").Value = 7239.474
$ws.Range("L77").Value = 10124.286
$ws.Range("M77").Value = -2871.474
$ws.Range("N77").Value = -18860.286
$ws.Range("H140").Value = 85714.22
$ws.Range("J140").Value = 85714.22
$ws.Range("L140").Value = 85714.22
$ws.Range("N140").Value = -96074.22

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 447
$ws.Range("J22").Value = 94
$ws.Range("L22").Value = 94
$ws.Range("N22").Value = -440

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14384.6
$ws.Range("I31").Value = 3258.7144
$ws.Range("K31").Value = 3258.7144
$ws.Range("M31").Value = -2963.7144
$ws.Range("H34").Value = 14384.6
$ws.Range("I34").Value = 3258.7144
$ws.Range("K34").Value = 3258.7144
$ws.Range("M34").Value = -3056.7144
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 203.4
$ws.Range("I10").Value = 197.71428
$ws.Range("J10").Value = 216.66667
$ws.Range("K10").Value = 593.14284
$ws.Range("L10").Value = 650.00001
$ws.Range("M10").Value = -454.14284
$ws.Range("N10").Value = -928.00001
$ws.Range("H34").Value = 456.8421
$ws.Range("J34").Value = 886.125
$ws.Range("L34").Value = 2658.375
$ws.Range("N34").Value = -2826.375
$ws.Range("H39").Value = 9774.625
$ws.Range("J39").Value = 9774.625
$ws.Range("L39").Value = 29323.875
$ws.Range("N39").Value = -29911.875
$ws.Range("H55").Value = 7353622.5
$ws.Range("I55").Value = 262.2
$ws.Range("J55").Value = 10417523
$ws.Range("K55").Value = 786.5999999999999
$ws.Range("L55").Value = 31252569
$ws.Range("M55").Value = -609.5999999999999
$ws.Range("N55").Value = -31252923
$ws.Range("H58").Value = 10500
$ws.Range("J58").Value = 15000
$ws.Range("L58").Value = 45000
$ws.Range("N58").Value = -45256
$ws.Range("H63").Value = 1900
$ws.Range("I63").Value = 1900
$ws.Range("K63").Value = 5700
$ws.Range("M63").Value = -4951
$ws.Range("H64").Value = 3682.6667
$ws.Range("I64").Value = 2298.25
$ws.Range("J64").Value = 4374.875
$ws.Range("K64").Value = 6894.75
$ws.Range("L64").Value = 13124.625
$ws.Range("M64").Value = -6624.75
$ws.Range("N64").Value = -13664.625
$ws.Range("H66").Value = 1900
$ws.Range("I66").Value = 1900
$ws.Range("K66").Value = 17100
$ws.Range("M66").Value = -13356
$ws.Range("H67").Value = 3682.6667
$ws.Range("I67").Value = 2298.25
$ws.Range("J67").Value = 4374.875
$ws.Range("K67").Value = 6894.75
$ws.Range("L67").Value = 13124.625
$ws.Range("M67").Value = -5958.75
$ws.Range("N67").Value = -14996.625
$ws.Range("H81").Value = 38471732
$ws.Range("I81").Value = 5271
$ws.Range("J81").Value = 71442984
$ws.Range("K81").Value = 15813
$ws.Range("L81").Value = 214328952
$ws.Range("M81").Value = -14690
$ws.Range("N81").Value = -214331198
$ws.Range("H82").Value = 8750
$ws.Range("J82").Value = 8750
$ws.Range("L82").Value = 26250
$ws.Range("N82").Value = -27062
$ws.Range("H84").Value = 38471732
$ws.Range("I84").Value = 5271
$ws.Range("J84").Value = 71442984
$ws.Range("K84").Value = 47439
$ws.Range("L84").Value = 642986856
$ws.Range("M84").Value = -41823
$ws.Range("N84").Value = -642998088
$ws.Range("H85").Value = 8750
$ws.Range("J85").Value = 8750
$ws.Range("L85").Value = 26250
$ws.Range("N85").Value = -29058
$ws.Range("H94").Value = 6304
$ws.Range("I94").Value = 4162
$ws.Range("J94").Value = 7375
$ws.Range("K94").Value = 12486
$ws.Range("L94").Value = 22125
$ws.Range("M94").Value = -11810
$ws.Range("N94").Value = -23477

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5409.7393
$ws.Range("J70").Value = 5317.1665
$ws.Range("L70").Value = 5317.1665
$ws.Range("N70").Value = -5857.1665
$ws.Range("H73").Value = 5409.7393
$ws.Range("J73").Value = 5317.1665
$ws.Range("L73").Value = 5317.1665
$ws.Range("N73").Value = -7189.1665
$ws.Range("H132").Value = 4460.1025
$ws.Range("I132").Value = 3646.3
$ws.Range("J132").Value = 7172.778
$ws.Range("K132").Value = 10938.9
$ws.Range("L132").Value = 21518.334
$ws.Range("M132").Value = -8408.900000000001
$ws.Range("N132").Value = -26578.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 110514.4
$ws.Range("I16").Value = 157499.28
$ws.Range("J16").Value = 883
$ws.Range("K16").Value = 157499.28
$ws.Range("L16").Value = 883
$ws.Range("M16").Value = -157329.28
$ws.Range("N16").Value = -1223
$ws.Range("H22").Value = 3528.7693
$ws.Range("I22").Value = 5666.3335
$ws.Range("J22").Value = 2887.5
$ws.Range("K22").Value = 5666.3335
$ws.Range("L22").Value = 2887.5
$ws.Range("M22").Value = -5371.3335
$ws.Range("N22").Value = -3477.5
$ws.Range("H27").Value = 3528.7693
$ws.Range("I27").Value = 5666.3335
$ws.Range("J27").Value = 2887.5
$ws.Range("K27").Value = 5666.3335
$ws.Range("L27").Value = 2887.5
$ws.Range("M27").Value = -5559.3335
$ws.Range("N27").Value = -3101.5
$ws.Range("H55").Value = 175.1579
$ws.Range("I55").Value = 171.16667
$ws.Range("J55").Value = 182
$ws.Range("K55").Value = 171.16667
$ws.Range("L55").Value = 182
$ws.Range("M55").Value = 1.833329999999989
$ws.Range("N55").Value = -528
$ws.Range("H60").Value = 28000
$ws.Range("J60").Value = 28000
$ws.Range("L60").Value = 28000
$ws.Range("N60").Value = -29018
$ws.Range("H68").Value = 1566.3334
$ws.Range("I68").Value = 699
$ws.Range("K68").Value = 699
$ws.Range("M68").Value = 50
$ws.Range("H71").Value = 1566.3334
$ws.Range("I71").Value = 699
$ws.Range("K71").Value = 3495
$ws.Range("M71").Value = 249
$ws.Range("H136").Value = 4826.15
$ws.Range("I136").Value = 4395.0713
$ws.Range("K136").Value = 13185.2139
$ws.Range("M136").Value = -10635.2139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 12027
$ws.Range("I43").Value = 12027
$ws.Range("K43").Value = 12027
$ws.Range("M43").Value = -11878
$ws.Range("H46").Value = 100000
$ws.Range("J46").Value = 100000
$ws.Range("L46").Value = 100000
$ws.Range("N46").Value = -100462
$ws.Range("H107").Value = 35715948
$ws.Range("I107").Value = 2291.75
$ws.Range("K107").Value = 6875.25
$ws.Range("M107").Value = -4955.25
$ws.Range("H122").Value = 1087.0869
$ws.Range("I122").Value = 974.7
$ws.Range("K122").Value = 2924.1
$ws.Range("M122").Value = -474.1000000000004
$ws.Range("H134").Value = 100000
$ws.Range("J134").Value = 100000
$ws.Range("L134").Value = 300000
$ws.Range("N134").Value = -305070

Write-Host "Applied all market-data cell updates."